# "Generate Report for Handback"
# The d01e54dd-f5f1-4163-97f9-d5a2625eda64.md file has been handed back.
# Update its Status on the zh-cn/de-de report sheets (and roll-up Status
# on the Overview sheet), and refresh the shared "Latest Handback
# DateTime" value that row 3 (8e14a4a7) and row 5 (d01e54dd) happened to
# have in common.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C5").Value = $newStatus
$wsZh.Range("H3").Value = "2016-03-22 06:24:04"
$wsZh.Range("H5").Value = "2016-03-22 06:24:04"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C5").Value = $newStatus
$wsDe.Range("H3").Value = "2016-03-22 06:24:10"
$wsDe.Range("H5").Value = "2016-03-22 06:24:10"

# Overview sheet (row 5 = d01e54dd file): zh-cn and de-de status columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B5").Value = $newStatus
$wsOverview.Range("C5").Value = $newStatus
